$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quarter labels for rows 2..39 (years 1987..2024, all Q4)
$startYear = 1987
$firstRow = 2
$lastRow = 39

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $year = $startYear + ($row - $firstRow)
    $label = "$($year)Q4"

    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $label
}

# Re-apply the same visual style (font/border/alignment/number-format) used
# in the header row (A1) to the newly text-ified date column so formatting
# matches exactly (this also clears the old date number format).
$headerCell = $ws.Range("A1")
$targetRange = $ws.Range("A$firstRow`:A$lastRow")
$headerCell.Copy()
$targetRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
